# Iteration og faseplan opdateret
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overordnet projektplan")

# Fill in the new "Aktiviteter" (F) / "Dato" (G) style columns for rows 6-11.
# Values pulled from the updated sharedStrings table.
$ws.Range("F6").Value = "System sekvens diagrammer"
$ws.Range("G6").Value = "2 timer"

$ws.Range("F7").Value = "Review af gruppe 2s artefakter"
$ws.Range("G7").Value = "0,5 timer"

$ws.Range("F8").Value = "Operations kontraker"
$ws.Range("G8").Value = "1 time"

$ws.Range("F9").Value = "rapport skrivning"
$ws.Range("F10").Value = "database model og normalformer"
$ws.Range("G9").Value = "3 timer"
$ws.Range("G10").Value = "3 timer"

$ws.Range("F11").Value = "Klassediagram med fokus på 3lags deling"
$ws.Range("G11").Value = "1 time"

# Apply the same cell style used by the neighboring cells (style index 20,
# i.e. plain style with applyFill/applyBorder) to the newly-populated cells
# in F9:G11, matching F6:G8 which already carried that style.
$ws.Range("F9:G11").Style = $ws.Range("F8").Style

# Update the frozen-pane top-left cell and the active selection to match
# the saved view state.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G11").Select()
